# Update the "Date" column (column B) for rows 2-14 with the new
# test-run timestamps captured for the VLink "Card Not Accepted" test
# cases, as described in the commit:
# "Added Test Cases and data for VLink Card Not Accepted Error."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = "Wed Feb 08 17:04:16 EST 2023"
$ws.Range("B3").Value  = "Wed Feb 08 17:05:04 EST 2023"
$ws.Range("B4").Value  = "Wed Feb 08 17:05:51 EST 2023"
$ws.Range("B5").Value  = "Wed Feb 08 17:06:37 EST 2023"
$ws.Range("B6").Value  = "Wed Feb 08 17:07:27 EST 2023"
$ws.Range("B7").Value  = "Wed Feb 08 17:08:15 EST 2023"
$ws.Range("B8").Value  = "Wed Feb 08 17:09:01 EST 2023"
$ws.Range("B9").Value  = "Wed Feb 08 17:09:49 EST 2023"
$ws.Range("B10").Value = "Wed Feb 08 17:10:37 EST 2023"
$ws.Range("B11").Value = "Wed Feb 08 17:11:25 EST 2023"
$ws.Range("B12").Value = "Wed Feb 08 17:12:14 EST 2023"
$ws.Range("B13").Value = "Wed Feb 08 17:13:02 EST 2023"
$ws.Range("B14").Value = "Wed Feb 08 17:13:49 EST 2023"

$wb.Save()
